$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row before current row 30, shifting existing
# rows 30-79 down to 31-80 (this matches the weekly price update reflected
# in the diff: dimension grows from A1:R79 to A1:R80).
$ws.Rows("30:30").Insert()

# Populate the newly inserted row 30 with the new weekly record.
$ws.Range("A30").Value = 7
$ws.Range("B30").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C30").Value = 'Ñuble'
$ws.Range("D30").Value = (Get-Date -Year 2022 -Month 8 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E30").Value = 16
$ws.Range("F30").Value = 100112040
$ws.Range("G30").Value = 'Cilantro'
$ws.Range("H30").Value = 'Sin especificar'
$ws.Range("I30").Value = 'Segunda'
$ws.Range("J30").Value = 100
$ws.Range("K30").Value = 600
$ws.Range("L30").Value = 600
$ws.Range("M30").Value = 600
$ws.Range("N30").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O30").Value = 'Provincia de Diguillín'
$ws.Range("P30").Value = 600
$ws.Range("Q30").Value = 1
$ws.Range("R30").Value = 'Hortaliza'
